# Apply the "only_quest_a" update: a new survey response (row 15) was added
# to the hidden query-result sheet "Video Game Narrative Measuremen", the
# free-text answer columns were shifted from AT:AW to AX:BA (to line up
# with the "Kolumna#" headers), and the table / defined name / dimension
# were extended to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Video Game Narrative Measuremen")

# ---------------------------------------------------------------------
# 1. Move the existing free-text answers (AT2:AW14) over to AX2:BA14 so
#    they line up under the "Kolumna1".."Kolumna4" headers instead of
#    the "How did you feel interacting with NPCs?" headers.
# ---------------------------------------------------------------------
$ws.Range("AT2:AW14").Cut($ws.Range("AX2:BA14"))

# ---------------------------------------------------------------------
# 2. Append the new response as row 15.
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "2024/06/05 2:18:16 PM EEST"
$ws.Range("C15").Value = "Male"
$ws.Range("D15").Value = 29
$ws.Range("E15").Value = "Sweden"
$ws.Range("F15").Value = "7-12"
$ws.Range("G15").Value = "10-15h"

$rowValues = @(4,3,1,1,3,4,4,2,2,3,1,3,3,4,4,4,2,4,1,4,4,3,4,4,4,4,2,2,4,2,1,2,4,1,2,3,4,1)
$col = 8
foreach ($val in $rowValues) {
    $ws.Cells.Item(15, $col).Value = $val
    $col++
}

$ws.Range("AX15").Value = "the npcs felt more alive but also harder to grasp since i had to ask questions"
$ws.Range("AY15").Value = "it was interesting trying to come up with good questions"
$ws.Range("AZ15").Value = "yes"

# ---------------------------------------------------------------------
# 3. Extend the query table so it covers the new row.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:BA15"))

# ---------------------------------------------------------------------
# 4. Update the hidden ExternalData_1 defined name to cover the new row.
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Video Game Narrative Measuremen!ExternalData_1") {
        $n.RefersTo = "='Video Game Narrative Measuremen'!`$A`$1:`$BA`$15"
    }
}

# ---------------------------------------------------------------------
# 5. Restore the selection on the data sheet to AA2:AS15.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("AA2:AS15").Select()

$wb.Save()
